# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-05-09 Friday", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "2025-05-10 Saturday", 2)

# Update the division problems in the single table. Addressing cells by
# (row, column) avoids ambiguity since "67÷3=" appears twice in the
# original table (row 9 col 1 and row 13 col 5) with different replacements.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "10÷9=" },
    @{ Row = 1;  Col = 2; Text = "81÷5=" },
    @{ Row = 1;  Col = 3; Text = "78÷3=" },
    @{ Row = 1;  Col = 4; Text = "62÷7=" },
    @{ Row = 1;  Col = 5; Text = "12÷8=" },

    @{ Row = 5;  Col = 1; Text = "65÷2=" },
    @{ Row = 5;  Col = 2; Text = "99÷3=" },
    @{ Row = 5;  Col = 3; Text = "23÷7=" },
    @{ Row = 5;  Col = 4; Text = "34÷2=" },
    @{ Row = 5;  Col = 5; Text = "48÷5=" },

    @{ Row = 9;  Col = 1; Text = "16÷4=" },
    @{ Row = 9;  Col = 2; Text = "70÷5=" },
    @{ Row = 9;  Col = 3; Text = "77÷8=" },
    @{ Row = 9;  Col = 4; Text = "48÷6=" },
    @{ Row = 9;  Col = 5; Text = "80÷5=" },

    @{ Row = 13; Col = 1; Text = "77÷6=" },
    @{ Row = 13; Col = 2; Text = "72÷9=" },
    @{ Row = 13; Col = 3; Text = "92÷6=" },
    @{ Row = 13; Col = 4; Text = "63÷9=" },
    @{ Row = 13; Col = 5; Text = "27÷6=" },

    @{ Row = 17; Col = 1; Text = "82÷8=" },
    @{ Row = 17; Col = 2; Text = "31÷6=" },
    @{ Row = 17; Col = 3; Text = "51÷5=" },
    @{ Row = 17; Col = 4; Text = "46÷6=" },
    @{ Row = 17; Col = 5; Text = "39÷2=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
